# Apply "Änderungen vom Labor (SW Test, ...)" edit
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AutoTest")

$checkMark = [char]0x2611   # ☑

function Set-Check($addr) {
    $target = $ws.Range($addr)
    $target.Value = $checkMark
    $target.HorizontalAlignment = -4108  # xlCenter
}

# Row 11: C,D,E,F,G,H,J get checkmarks (I11 and K11 stay untouched/empty)
foreach ($col in @("C","D","E","F","G","H","J")) {
    Set-Check ($col + "11")
}

# Row 13: D,E,F,H,I,J get checkmarks (C13,G13 already had them; K13 untouched)
foreach ($col in @("D","E","F","H","I","J")) {
    Set-Check ($col + "13")
}

# Rows 14-17: D,E,F,H get checkmarks (C,G already had them; I,J,K untouched)
foreach ($r in 14..17) {
    foreach ($col in @("D","E","F","H")) {
        Set-Check ($col + $r)
    }
}

# Row 18: H gets checkmark (C18,G18 already had them)
Set-Check "H18"

# Update the active cell selection to I22
$ws.Range("I22").Select()
